# Perft - Character section
# Adds a new results column (I) for the "After fixing character section" run,
# mirroring the layout already used by the other run columns (B..G),
# and introduces a narrow spacer column (H) between the existing data
# and the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow "spacer" column between the existing data (B:G) and the new
# run column (I). The target width is 10.42578125 characters; this
# runtime quantizes custom widths, so we pick the closest obtainable
# ColumnWidth value.
$ws.Columns.Item(8).ColumnWidth = 9.6

# Header for the new run.
$ws.Range("I3").Value = "After fixing character section"

# Block 1
$ws.Range("I5").Value = "RunPerft1"
$ws.Range("I6").Value = "Voices,Seconds,TimeMillis"
$ws.Range("I7").Value = "1,20,528"

# Block 2
$ws.Range("I8").Value = "RunPerft1"
$ws.Range("I9").Value = "Voices,Seconds,TimeMillis"
$ws.Range("I10").Value = "20,2,984"

# Block 3
$ws.Range("I11").Value = "RunPerft1"
$ws.Range("I12").Value = "Voices,Seconds,TimeMillis"
$ws.Range("I13").Value = "1,20,519"

# Block 4
$ws.Range("I14").Value = "RunPerft1"
$ws.Range("I15").Value = "Voices,Seconds,TimeMillis"
$ws.Range("I16").Value = "20,2,989"

# Block 5
$ws.Range("I17").Value = "RunPerft1"
$ws.Range("I18").Value = "Voices,Seconds,TimeMillis"
$ws.Range("I19").Value = "1,20,516"

# Block 6
$ws.Range("I20").Value = "RunPerft1"
$ws.Range("I21").Value = "Voices,Seconds,TimeMillis"
$ws.Range("I22").Value = "20,2,985"

# Scroll the view toward the new column and move the selection to match
# the author's final cursor position.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I16").Select()
